# + show entries datatables & perbaikan impor periode & penamaan file export
#
# This template workbook ("Template Import User Alumni") carries example
# values that illustrate the "periode" (academic period) column format.
# The commit fixes those example/sample values and tidies up a leftover,
# unused cell style.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1 - example data row
$ws2 = $wb.Worksheets.Item(2)   # Sheet2 - instructions / legend

# ---------------------------------------------------------------------
# 1) "periode" example values
#    Sheet1!E1 held the old sample "2022/2023 Genap" -> replace with the
#    new sample period "2024/2025" (perbaikan impor periode).
# ---------------------------------------------------------------------
$ws1.Range("E1").Value = "2024/2025"

# ---------------------------------------------------------------------
# 2) Sheet2!E2 documents the expected "periode" format and embeds an
#    example value inline in its help text; the example changes from
#    "2022/2024 Genap" to just "2022/2024".
# ---------------------------------------------------------------------
$ws2.Range("E2").Value = "format penulisan yang berlaku pada kolom periode mengunakan tahun/tahun jenis periode (contoh 2022/2024)"

# ---------------------------------------------------------------------
# 3) Sheet2!D3 was carrying a stray/duplicate cell style (an empty
#    alignment definition) instead of the plain style shared by its
#    row neighbours (C3, E3, F3). Re-apply the neighbour's formatting
#    so the redundant style is no longer used.
# ---------------------------------------------------------------------
[void]$ws2.Range("C3").Copy()
[void]$ws2.Range("D3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Refresh the in-workbook selection/navigation state left over from
#    editing, matching where the author ended up when they saved.
# ---------------------------------------------------------------------
$ws2.Activate()
[void]$ws2.Range("E17").Select()

$ws1.Activate()
[void]$ws1.Range("F9").Select()
